# fix: contrucción grafico de resultados
#
# The "Recursos" sheet listed volumes for "BRA caña" (row 8) in the wrong
# unit (tonnes instead of kg), which threw off the results chart built from
# that data. Correct the four affected cells by the missing factor of 1000,
# then restore each sheet's last-used selection.

$wb = $excel.ActiveWorkbook

$wsRecursos = $wb.Worksheets.Item("Recursos")
$wsDatos    = $wb.Worksheets.Item("Datos")

# --- data fix: row 8 ("BRA caña") was off by a factor of 1000 -------------
$wsRecursos.Range("C8").Value = 1406000
$wsRecursos.Range("D8").Value = 2813000
$wsRecursos.Range("E8").Value = 4219000
$wsRecursos.Range("G8").Value = 14064000

# --- restore the selection left on the "Datos" sheet -----------------------
[void]$wsDatos.Range("D18").Select()

# --- restore "Recursos" as the active sheet/selection (must be last so it
#     stays the active tab) --------------------------------------------------
[void]$wsRecursos.Activate()
[void]$wsRecursos.Range("G11").Select()
